$p = $ppt.ActivePresentation
$s = $p.Slides.Add(8, 12)
Write-Host "Added slide, count is now $($p.Slides.Count)"
